# ---------------------------------------------------------------------------
# Adds a new "PO Forecast" worksheet (ds / PO_Forecast / yhat_lower /
# yhat_upper) after the existing sheets, and renames the "Requested
# quantity" headers on the two existing sheets to Weekly_PO_Qty /
# Monthly_PO_Qty respectively.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1. Rename the "Requested quantity" header on "Weekly Quantity" --------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

# --- 2. Rename the "Requested quantity" header on "Monthly Trend" ---------
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 3. Add the new "PO Forecast" worksheet at the end of the workbook ----
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Match the outline / page-setup properties used by the other sheets.
$wsForecast.Outline.SummaryRow = 1
$wsForecast.Outline.SummaryColumn = 1
$wsForecast.PageSetup.LeftMargin = 54
$wsForecast.PageSetup.RightMargin = 54
$wsForecast.PageSetup.TopMargin = 72
$wsForecast.PageSetup.BottomMargin = 72
$wsForecast.PageSetup.HeaderMargin = 36
$wsForecast.PageSetup.FooterMargin = 36

# Header row.
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Copy the bold/centered header style from an existing sheet's header row.
$wsWeekly.Range("A1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122) # xlPasteFormats

# Data rows (82 rows of forecast data), written in one shot via a 2-D array.
$rowCount = 82
$arr = New-Object 'object[,]' $rowCount,4

$arr[0,0]=44941.99999999999; $arr[0,1]=75; $arr[0,2]=-37.15906299465209; $arr[0,3]=189.59216873859
$arr[1,0]=44948.99999999999; $arr[1,1]=75; $arr[1,2]=-40.20891278447392; $arr[1,3]=184.1286742061823
$arr[2,0]=44969.99999999999; $arr[2,1]=75; $arr[2,2]=-39.60442762956134; $arr[2,3]=194.4967218824556
$arr[3,0]=44976.99999999999; $arr[3,1]=75; $arr[3,2]=-35.17360466017416; $arr[3,3]=187.3366157872235
$arr[4,0]=44983.99999999999; $arr[4,1]=75; $arr[4,2]=-37.78144490955377; $arr[4,3]=196.482177846407
$arr[5,0]=44990.99999999999; $arr[5,1]=76; $arr[5,2]=-27.83208154300777; $arr[5,3]=182.2944508119008
$arr[6,0]=44997.99999999999; $arr[6,1]=76; $arr[6,2]=-43.03515247082616; $arr[6,3]=193.6212459503906
$arr[7,0]=45004.99999999999; $arr[7,1]=76; $arr[7,2]=-43.58375433746596; $arr[7,3]=188.0021321332204
$arr[8,0]=45011.99999999999; $arr[8,1]=76; $arr[8,2]=-43.95564888273719; $arr[8,3]=188.1223192860888
$arr[9,0]=45025.99999999999; $arr[9,1]=76; $arr[9,2]=-26.12745331639386; $arr[9,3]=188.3123875005852
$arr[10,0]=45032.99999999999; $arr[10,1]=76; $arr[10,2]=-34.33682354674438; $arr[10,3]=192.3532600411637
$arr[11,0]=45039.99999999999; $arr[11,1]=76; $arr[11,2]=-36.35869881970251; $arr[11,3]=183.7722036760842
$arr[12,0]=45053.99999999999; $arr[12,1]=77; $arr[12,2]=-39.72042010502889; $arr[12,3]=198.6625174764845
$arr[13,0]=45060.99999999999; $arr[13,1]=77; $arr[13,2]=-35.91786697927918; $arr[13,3]=190.2897080093356
$arr[14,0]=45067.99999999999; $arr[14,1]=77; $arr[14,2]=-32.29520371916414; $arr[14,3]=197.3456319851228
$arr[15,0]=45074.99999999999; $arr[15,1]=77; $arr[15,2]=-36.19159352395357; $arr[15,3]=194.2578682705242
$arr[16,0]=45088.99999999999; $arr[16,1]=77; $arr[16,2]=-29.12962287704124; $arr[16,3]=189.6614344057303
$arr[17,0]=45095.99999999999; $arr[17,1]=77; $arr[17,2]=-32.52790480156173; $arr[17,3]=189.7701126514765
$arr[18,0]=45102.99999999999; $arr[18,1]=77; $arr[18,2]=-39.14372220340869; $arr[18,3]=186.8402345350079
$arr[19,0]=45109.99999999999; $arr[19,1]=78; $arr[19,2]=-40.11962122013392; $arr[19,3]=189.0347167832686
$arr[20,0]=45116.99999999999; $arr[20,1]=78; $arr[20,2]=-39.39778035418804; $arr[20,3]=184.9679454295316
$arr[21,0]=45123.99999999999; $arr[21,1]=78; $arr[21,2]=-40.655029248471; $arr[21,3]=187.0739663998472
$arr[22,0]=45130.99999999999; $arr[22,1]=78; $arr[22,2]=-30.5849608986286; $arr[22,3]=186.3488818936045
$arr[23,0]=45137.99999999999; $arr[23,1]=78; $arr[23,2]=-36.43015174784971; $arr[23,3]=199.4105099613044
$arr[24,0]=45144.99999999999; $arr[24,1]=78; $arr[24,2]=-41.07317265808352; $arr[24,3]=199.9416737176444
$arr[25,0]=45151.99999999999; $arr[25,1]=78; $arr[25,2]=-41.00243919500274; $arr[25,3]=193.0219715910637
$arr[26,0]=45158.99999999999; $arr[26,1]=78; $arr[26,2]=-34.33091102088202; $arr[26,3]=192.3704214506515
$arr[27,0]=45165.99999999999; $arr[27,1]=78; $arr[27,2]=-21.61198691120789; $arr[27,3]=193.4344953653406
$arr[28,0]=45172.99999999999; $arr[28,1]=79; $arr[28,2]=-37.93952803803821; $arr[28,3]=183.3656955771145
$arr[29,0]=45179.99999999999; $arr[29,1]=79; $arr[29,2]=-38.82280200349515; $arr[29,3]=195.4972082799594
$arr[30,0]=45186.99999999999; $arr[30,1]=79; $arr[30,2]=-39.51405560408833; $arr[30,3]=195.4225234399433
$arr[31,0]=45193.99999999999; $arr[31,1]=79; $arr[31,2]=-29.28995023770346; $arr[31,3]=187.8982064563286
$arr[32,0]=45200.99999999999; $arr[32,1]=79; $arr[32,2]=-35.27795732699293; $arr[32,3]=194.6245831187207
$arr[33,0]=45207.99999999999; $arr[33,1]=79; $arr[33,2]=-27.75496846260449; $arr[33,3]=182.0919569479525
$arr[34,0]=45214.99999999999; $arr[34,1]=79; $arr[34,2]=-43.33270954017764; $arr[34,3]=184.5513970934623
$arr[35,0]=45221.99999999999; $arr[35,1]=79; $arr[35,2]=-26.70374714977017; $arr[35,3]=185.1661601292481
$arr[36,0]=45228.99999999999; $arr[36,1]=79; $arr[36,2]=-33.34891316590097; $arr[36,3]=196.0669107331753
$arr[37,0]=45235.99999999999; $arr[37,1]=80; $arr[37,2]=-25.66628058133994; $arr[37,3]=188.1505527002854
$arr[38,0]=45242.99999999999; $arr[38,1]=80; $arr[38,2]=-36.47697792573; $arr[38,3]=193.0425286781153
$arr[39,0]=45249.99999999999; $arr[39,1]=80; $arr[39,2]=-32.47220340339918; $arr[39,3]=193.3075664407905
$arr[40,0]=45256.99999999999; $arr[40,1]=80; $arr[40,2]=-37.3624414667025; $arr[40,3]=195.2332850915811
$arr[41,0]=45270.99999999999; $arr[41,1]=80; $arr[41,2]=-29.53004822438182; $arr[41,3]=196.2167181291339
$arr[42,0]=45277.99999999999; $arr[42,1]=80; $arr[42,2]=-22.38765558317026; $arr[42,3]=200.5132006177361
$arr[43,0]=45298.99999999999; $arr[43,1]=81; $arr[43,2]=-30.18608140861252; $arr[43,3]=185.6857557741064
$arr[44,0]=45305.99999999999; $arr[44,1]=81; $arr[44,2]=-29.85605578665161; $arr[44,3]=198.0922579228793
$arr[45,0]=45312.99999999999; $arr[45,1]=81; $arr[45,2]=-32.68549839110852; $arr[45,3]=192.2270255953368
$arr[46,0]=45319.99999999999; $arr[46,1]=81; $arr[46,2]=-26.91478569613734; $arr[46,3]=197.0849978652881
$arr[47,0]=45326.99999999999; $arr[47,1]=81; $arr[47,2]=-32.48488467992205; $arr[47,3]=198.2397837892332
$arr[48,0]=45333.99999999999; $arr[48,1]=81; $arr[48,2]=-29.50305515867789; $arr[48,3]=197.6664251496616
$arr[49,0]=45340.99999999999; $arr[49,1]=81; $arr[49,2]=-34.40282325433546; $arr[49,3]=201.6278745912924
$arr[50,0]=45347.99999999999; $arr[50,1]=81; $arr[50,2]=-23.56401821281846; $arr[50,3]=194.6274374516667
$arr[51,0]=45354.99999999999; $arr[51,1]=81; $arr[51,2]=-33.57732486334828; $arr[51,3]=198.6098484291531
$arr[52,0]=45361.99999999999; $arr[52,1]=82; $arr[52,2]=-28.10693330259602; $arr[52,3]=200.0272887286415
$arr[53,0]=45368.99999999999; $arr[53,1]=82; $arr[53,2]=-25.86453997161984; $arr[53,3]=198.2846531775762
$arr[54,0]=45375.99999999999; $arr[54,1]=82; $arr[54,2]=-29.20677984538684; $arr[54,3]=208.9495224181285
$arr[55,0]=45382.99999999999; $arr[55,1]=82; $arr[55,2]=-39.47250113765695; $arr[55,3]=201.5844474267664
$arr[56,0]=45389.99999999999; $arr[56,1]=82; $arr[56,2]=-30.62868859278155; $arr[56,3]=196.2332217229595
$arr[57,0]=45396.99999999999; $arr[57,1]=82; $arr[57,2]=-33.0955866084788; $arr[57,3]=191.4410185246331
$arr[58,0]=45403.99999999999; $arr[58,1]=82; $arr[58,2]=-22.16302629529122; $arr[58,3]=197.1865744293893
$arr[59,0]=45410.99999999999; $arr[59,1]=82; $arr[59,2]=-27.88258861969993; $arr[59,3]=197.4127632207089
$arr[60,0]=45424.99999999999; $arr[60,1]=83; $arr[60,2]=-28.78549185283353; $arr[60,3]=192.4976896151888
$arr[61,0]=45487.99999999999; $arr[61,1]=84; $arr[61,2]=-24.2553038920175; $arr[61,3]=193.3900967802895
$arr[62,0]=45494.99999999999; $arr[62,1]=84; $arr[62,2]=-31.97082868342194; $arr[62,3]=196.796684707995
$arr[63,0]=45501.99999999999; $arr[63,1]=84; $arr[63,2]=-31.70815046647904; $arr[63,3]=194.247927407999
$arr[64,0]=45515.99999999999; $arr[64,1]=84; $arr[64,2]=-32.47792347031545; $arr[64,3]=203.6695536139865
$arr[65,0]=45529.99999999999; $arr[65,1]=84; $arr[65,2]=-21.07266674374726; $arr[65,3]=203.2517097607161
$arr[66,0]=45536.99999999999; $arr[66,1]=84; $arr[66,2]=-30.33951338372975; $arr[66,3]=196.1366277903282
$arr[67,0]=45543.99999999999; $arr[67,1]=85; $arr[67,2]=-28.86872441443037; $arr[67,3]=198.0209994667541
$arr[68,0]=45550.99999999999; $arr[68,1]=85; $arr[68,2]=-31.40941463233712; $arr[68,3]=199.737033539323
$arr[69,0]=45557.99999999999; $arr[69,1]=85; $arr[69,2]=-29.94504807473989; $arr[69,3]=190.051078888518
$arr[70,0]=45578.99999999999; $arr[70,1]=85; $arr[70,2]=-30.79035913376808; $arr[70,3]=199.5023061291029
$arr[71,0]=45592.99999999999; $arr[71,1]=85; $arr[71,2]=-29.1993030666475; $arr[71,3]=190.7828921081615
$arr[72,0]=45599.99999999999; $arr[72,1]=85; $arr[72,2]=-24.27112490492829; $arr[72,3]=194.5489516358918
$arr[73,0]=45613.99999999999; $arr[73,1]=86; $arr[73,2]=-24.90277966297209; $arr[73,3]=200.9734688664034
$arr[74,0]=45620.99999999999; $arr[74,1]=86; $arr[74,2]=-22.55762260775614; $arr[74,3]=201.2615384874247
$arr[75,0]=45627.99999999999; $arr[75,1]=86; $arr[75,2]=-25.86642680672004; $arr[75,3]=195.1952149339868
$arr[76,0]=45634.99999999999; $arr[76,1]=86; $arr[76,2]=-20.80081276501289; $arr[76,3]=193.233626987226
$arr[77,0]=45641.99999999999; $arr[77,1]=86; $arr[77,2]=-19.00584546884331; $arr[77,3]=201.1400373183608
$arr[78,0]=45648.99999999999; $arr[78,1]=86; $arr[78,2]=-17.47097128148377; $arr[78,3]=202.5492666437016
$arr[79,0]=45655.99999999999; $arr[79,1]=86; $arr[79,2]=-24.97551244851524; $arr[79,3]=195.920458704838
$arr[80,0]=45662.99999999999; $arr[80,1]=86; $arr[80,2]=-29.21966965735516; $arr[80,3]=204.4870207567477
$arr[81,0]=45669.99999999999; $arr[81,1]=87; $arr[81,2]=-26.2153279952077; $arr[81,3]=201.265739935869

$wsForecast.Range("A2:D83").Value = $arr

# Copy the date/time number-format style used on the other sheets' date
# column onto the new sheet's "ds" column (A2:A83).
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A83").PasteSpecial(-4122) # xlPasteFormats

$excel.CutCopyMode = $false

# Keep the first sheet as the active / selected tab, same as before the edit.
$wb.Worksheets.Item(1).Activate()

Write-Output "PO Forecast sheet added"
